# ---------------------------------------------------------------------------
# "Fruta / hortaliza, semanal" weekly update.
#
# A new week of price observations (the "Primera"/"Segunda" pair now dated
# 2022-12-22) is published at the top of the data block (rows 176-177, which
# only get a refreshed date). Every older weekly pair that used to occupy rows
# 178-256 shifts two rows down (i.e. row N now holds what row N-2 held before),
# and the oldest pair that used to be the final rows 255-256 is preserved by
# appending it as two brand-new rows, 257-258, growing the sheet from
# A1:R256 to A1:R258.
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = $ws.Cells.Item(256, 4).NumberFormat

$ws.Cells.Item(176, 4).Value = 44917
$ws.Cells.Item(176, 9).Value = "Primera"
$ws.Cells.Item(176, 10).Value = 200
$ws.Cells.Item(176, 11).Value = 600
$ws.Cells.Item(176, 12).Value = 700
$ws.Cells.Item(176, 13).Value = 650
$ws.Cells.Item(176, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(176, 15).Value = "Región de Ñuble"
$ws.Cells.Item(176, 16).Value = 650
$ws.Cells.Item(176, 17).Value = 1

$ws.Cells.Item(177, 4).Value = 44917
$ws.Cells.Item(177, 9).Value = "Segunda"
$ws.Cells.Item(177, 10).Value = 100
$ws.Cells.Item(177, 11).Value = 500
$ws.Cells.Item(177, 12).Value = 500
$ws.Cells.Item(177, 13).Value = 500
$ws.Cells.Item(177, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(177, 15).Value = "Región de Ñuble"
$ws.Cells.Item(177, 16).Value = 500
$ws.Cells.Item(177, 17).Value = 1

$ws.Cells.Item(178, 4).Value = 44336
$ws.Cells.Item(178, 9).Value = "Primera"
$ws.Cells.Item(178, 10).Value = 200
$ws.Cells.Item(178, 11).Value = 600
$ws.Cells.Item(178, 12).Value = 700
$ws.Cells.Item(178, 13).Value = 650
$ws.Cells.Item(178, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(178, 15).Value = "Región de Ñuble"
$ws.Cells.Item(178, 16).Value = 650
$ws.Cells.Item(178, 17).Value = 1

$ws.Cells.Item(179, 4).Value = 44336
$ws.Cells.Item(179, 9).Value = "Segunda"
$ws.Cells.Item(179, 10).Value = 100
$ws.Cells.Item(179, 11).Value = 500
$ws.Cells.Item(179, 12).Value = 500
$ws.Cells.Item(179, 13).Value = 500
$ws.Cells.Item(179, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(179, 15).Value = "Región de Ñuble"
$ws.Cells.Item(179, 16).Value = 500
$ws.Cells.Item(179, 17).Value = 1

$ws.Cells.Item(180, 4).Value = 44775
$ws.Cells.Item(180, 9).Value = "Primera"
$ws.Cells.Item(180, 10).Value = 150
$ws.Cells.Item(180, 11).Value = 8500
$ws.Cells.Item(180, 12).Value = 9000
$ws.Cells.Item(180, 13).Value = 8733
$ws.Cells.Item(180, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(180, 15).Value = "Región Metropolitana"
$ws.Cells.Item(180, 16).Value = 243
$ws.Cells.Item(180, 17).Value = 36

$ws.Cells.Item(181, 4).Value = 44727
$ws.Cells.Item(181, 9).Value = "Primera"
$ws.Cells.Item(181, 10).Value = 170
$ws.Cells.Item(181, 11).Value = 5000
$ws.Cells.Item(181, 12).Value = 5500
$ws.Cells.Item(181, 13).Value = 5265
$ws.Cells.Item(181, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(181, 15).Value = "Región Metropolitana"
$ws.Cells.Item(181, 16).Value = 146
$ws.Cells.Item(181, 17).Value = 36

$ws.Cells.Item(182, 4).Value = 44750
$ws.Cells.Item(182, 9).Value = "Primera"
$ws.Cells.Item(182, 10).Value = 140
$ws.Cells.Item(182, 11).Value = 9500
$ws.Cells.Item(182, 12).Value = 10000
$ws.Cells.Item(182, 13).Value = 9714
$ws.Cells.Item(182, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(182, 15).Value = "Región Metropolitana"
$ws.Cells.Item(182, 16).Value = 270
$ws.Cells.Item(182, 17).Value = 36

$ws.Cells.Item(183, 4).Value = 44769
$ws.Cells.Item(183, 9).Value = "Primera"
$ws.Cells.Item(183, 10).Value = 200
$ws.Cells.Item(183, 11).Value = 700
$ws.Cells.Item(183, 12).Value = 800
$ws.Cells.Item(183, 13).Value = 750
$ws.Cells.Item(183, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(183, 15).Value = "Región de Ñuble"
$ws.Cells.Item(183, 16).Value = 750
$ws.Cells.Item(183, 17).Value = 1

$ws.Cells.Item(184, 4).Value = 44769
$ws.Cells.Item(184, 9).Value = "Segunda"
$ws.Cells.Item(184, 10).Value = 100
$ws.Cells.Item(184, 11).Value = 600
$ws.Cells.Item(184, 12).Value = 600
$ws.Cells.Item(184, 13).Value = 600
$ws.Cells.Item(184, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(184, 15).Value = "Región de Ñuble"
$ws.Cells.Item(184, 16).Value = 600
$ws.Cells.Item(184, 17).Value = 1

$ws.Cells.Item(185, 4).Value = 44705
$ws.Cells.Item(185, 9).Value = "Primera"
$ws.Cells.Item(185, 10).Value = 200
$ws.Cells.Item(185, 11).Value = 600
$ws.Cells.Item(185, 12).Value = 700
$ws.Cells.Item(185, 13).Value = 650
$ws.Cells.Item(185, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(185, 15).Value = "Región de Ñuble"
$ws.Cells.Item(185, 16).Value = 650
$ws.Cells.Item(185, 17).Value = 1

$ws.Cells.Item(186, 4).Value = 44705
$ws.Cells.Item(186, 9).Value = "Segunda"
$ws.Cells.Item(186, 10).Value = 100
$ws.Cells.Item(186, 11).Value = 500
$ws.Cells.Item(186, 12).Value = 500
$ws.Cells.Item(186, 13).Value = 500
$ws.Cells.Item(186, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(186, 15).Value = "Región de Ñuble"
$ws.Cells.Item(186, 16).Value = 500
$ws.Cells.Item(186, 17).Value = 1

$ws.Cells.Item(187, 4).Value = 44595
$ws.Cells.Item(187, 9).Value = "Primera"
$ws.Cells.Item(187, 10).Value = 200
$ws.Cells.Item(187, 11).Value = 600
$ws.Cells.Item(187, 12).Value = 700
$ws.Cells.Item(187, 13).Value = 650
$ws.Cells.Item(187, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(187, 15).Value = "Región de Ñuble"
$ws.Cells.Item(187, 16).Value = 650
$ws.Cells.Item(187, 17).Value = 1

$ws.Cells.Item(188, 4).Value = 44595
$ws.Cells.Item(188, 9).Value = "Segunda"
$ws.Cells.Item(188, 10).Value = 100
$ws.Cells.Item(188, 11).Value = 500
$ws.Cells.Item(188, 12).Value = 500
$ws.Cells.Item(188, 13).Value = 500
$ws.Cells.Item(188, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(188, 15).Value = "Región de Ñuble"
$ws.Cells.Item(188, 16).Value = 500
$ws.Cells.Item(188, 17).Value = 1

$ws.Cells.Item(189, 4).Value = 44442
$ws.Cells.Item(189, 9).Value = "Primera"
$ws.Cells.Item(189, 10).Value = 200
$ws.Cells.Item(189, 11).Value = 600
$ws.Cells.Item(189, 12).Value = 700
$ws.Cells.Item(189, 13).Value = 650
$ws.Cells.Item(189, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(189, 15).Value = "Región de Ñuble"
$ws.Cells.Item(189, 16).Value = 650
$ws.Cells.Item(189, 17).Value = 1

$ws.Cells.Item(190, 4).Value = 44442
$ws.Cells.Item(190, 9).Value = "Segunda"
$ws.Cells.Item(190, 10).Value = 100
$ws.Cells.Item(190, 11).Value = 500
$ws.Cells.Item(190, 12).Value = 500
$ws.Cells.Item(190, 13).Value = 500
$ws.Cells.Item(190, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(190, 15).Value = "Región de Ñuble"
$ws.Cells.Item(190, 16).Value = 500
$ws.Cells.Item(190, 17).Value = 1

$ws.Cells.Item(191, 4).Value = 44657
$ws.Cells.Item(191, 9).Value = "Primera"
$ws.Cells.Item(191, 10).Value = 180
$ws.Cells.Item(191, 11).Value = 6000
$ws.Cells.Item(191, 12).Value = 6500
$ws.Cells.Item(191, 13).Value = 6222
$ws.Cells.Item(191, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(191, 15).Value = "Región Metropolitana"
$ws.Cells.Item(191, 16).Value = 173
$ws.Cells.Item(191, 17).Value = 36

$ws.Cells.Item(192, 4).Value = 44685
$ws.Cells.Item(192, 9).Value = "Primera"
$ws.Cells.Item(192, 10).Value = 170
$ws.Cells.Item(192, 11).Value = 5500
$ws.Cells.Item(192, 12).Value = 6000
$ws.Cells.Item(192, 13).Value = 5765
$ws.Cells.Item(192, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(192, 15).Value = "Región Metropolitana"
$ws.Cells.Item(192, 16).Value = 160
$ws.Cells.Item(192, 17).Value = 36

$ws.Cells.Item(193, 4).Value = 44460
$ws.Cells.Item(193, 9).Value = "Primera"
$ws.Cells.Item(193, 10).Value = 200
$ws.Cells.Item(193, 11).Value = 600
$ws.Cells.Item(193, 12).Value = 700
$ws.Cells.Item(193, 13).Value = 650
$ws.Cells.Item(193, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(193, 15).Value = "Región de Ñuble"
$ws.Cells.Item(193, 16).Value = 650
$ws.Cells.Item(193, 17).Value = 1

$ws.Cells.Item(194, 4).Value = 44460
$ws.Cells.Item(194, 9).Value = "Segunda"
$ws.Cells.Item(194, 10).Value = 100
$ws.Cells.Item(194, 11).Value = 500
$ws.Cells.Item(194, 12).Value = 500
$ws.Cells.Item(194, 13).Value = 500
$ws.Cells.Item(194, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(194, 15).Value = "Región de Ñuble"
$ws.Cells.Item(194, 16).Value = 500
$ws.Cells.Item(194, 17).Value = 1

$ws.Cells.Item(195, 4).Value = 44777
$ws.Cells.Item(195, 9).Value = "Primera"
$ws.Cells.Item(195, 10).Value = 120
$ws.Cells.Item(195, 11).Value = 8000
$ws.Cells.Item(195, 12).Value = 8500
$ws.Cells.Item(195, 13).Value = 8292
$ws.Cells.Item(195, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(195, 15).Value = "Región Metropolitana"
$ws.Cells.Item(195, 16).Value = 230
$ws.Cells.Item(195, 17).Value = 36

$ws.Cells.Item(196, 4).Value = 44194
$ws.Cells.Item(196, 9).Value = "Primera"
$ws.Cells.Item(196, 10).Value = 200
$ws.Cells.Item(196, 11).Value = 600
$ws.Cells.Item(196, 12).Value = 700
$ws.Cells.Item(196, 13).Value = 650
$ws.Cells.Item(196, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(196, 15).Value = "Región de Ñuble"
$ws.Cells.Item(196, 16).Value = 650
$ws.Cells.Item(196, 17).Value = 1

$ws.Cells.Item(197, 4).Value = 44194
$ws.Cells.Item(197, 9).Value = "Segunda"
$ws.Cells.Item(197, 10).Value = 100
$ws.Cells.Item(197, 11).Value = 500
$ws.Cells.Item(197, 12).Value = 500
$ws.Cells.Item(197, 13).Value = 500
$ws.Cells.Item(197, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(197, 15).Value = "Región de Ñuble"
$ws.Cells.Item(197, 16).Value = 500
$ws.Cells.Item(197, 17).Value = 1

$ws.Cells.Item(198, 4).Value = 44322
$ws.Cells.Item(198, 9).Value = "Primera"
$ws.Cells.Item(198, 10).Value = 200
$ws.Cells.Item(198, 11).Value = 600
$ws.Cells.Item(198, 12).Value = 700
$ws.Cells.Item(198, 13).Value = 650
$ws.Cells.Item(198, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(198, 15).Value = "Región de Ñuble"
$ws.Cells.Item(198, 16).Value = 650
$ws.Cells.Item(198, 17).Value = 1

$ws.Cells.Item(199, 4).Value = 44322
$ws.Cells.Item(199, 9).Value = "Segunda"
$ws.Cells.Item(199, 10).Value = 100
$ws.Cells.Item(199, 11).Value = 500
$ws.Cells.Item(199, 12).Value = 500
$ws.Cells.Item(199, 13).Value = 500
$ws.Cells.Item(199, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(199, 15).Value = "Región de Ñuble"
$ws.Cells.Item(199, 16).Value = 500
$ws.Cells.Item(199, 17).Value = 1

$ws.Cells.Item(200, 4).Value = 44344
$ws.Cells.Item(200, 9).Value = "Primera"
$ws.Cells.Item(200, 10).Value = 200
$ws.Cells.Item(200, 11).Value = 600
$ws.Cells.Item(200, 12).Value = 700
$ws.Cells.Item(200, 13).Value = 650
$ws.Cells.Item(200, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(200, 15).Value = "Región de Ñuble"
$ws.Cells.Item(200, 16).Value = 650
$ws.Cells.Item(200, 17).Value = 1

$ws.Cells.Item(201, 4).Value = 44344
$ws.Cells.Item(201, 9).Value = "Segunda"
$ws.Cells.Item(201, 10).Value = 100
$ws.Cells.Item(201, 11).Value = 500
$ws.Cells.Item(201, 12).Value = 500
$ws.Cells.Item(201, 13).Value = 500
$ws.Cells.Item(201, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(201, 15).Value = "Región de Ñuble"
$ws.Cells.Item(201, 16).Value = 500
$ws.Cells.Item(201, 17).Value = 1

$ws.Cells.Item(202, 4).Value = 44761
$ws.Cells.Item(202, 9).Value = "Primera"
$ws.Cells.Item(202, 10).Value = 30
$ws.Cells.Item(202, 11).Value = 16000
$ws.Cells.Item(202, 12).Value = 17000
$ws.Cells.Item(202, 13).Value = 16333
$ws.Cells.Item(202, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(202, 15).Value = "Región Metropolitana"
$ws.Cells.Item(202, 16).Value = 454
$ws.Cells.Item(202, 17).Value = 36

$ws.Cells.Item(203, 4).Value = 44882
$ws.Cells.Item(203, 9).Value = "Primera"
$ws.Cells.Item(203, 10).Value = 200
$ws.Cells.Item(203, 11).Value = 700
$ws.Cells.Item(203, 12).Value = 800
$ws.Cells.Item(203, 13).Value = 750
$ws.Cells.Item(203, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(203, 15).Value = "Región de Ñuble"
$ws.Cells.Item(203, 16).Value = 750
$ws.Cells.Item(203, 17).Value = 1

$ws.Cells.Item(204, 4).Value = 44882
$ws.Cells.Item(204, 9).Value = "Segunda"
$ws.Cells.Item(204, 10).Value = 100
$ws.Cells.Item(204, 11).Value = 600
$ws.Cells.Item(204, 12).Value = 600
$ws.Cells.Item(204, 13).Value = 600
$ws.Cells.Item(204, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(204, 15).Value = "Región de Ñuble"
$ws.Cells.Item(204, 16).Value = 600
$ws.Cells.Item(204, 17).Value = 1

$ws.Cells.Item(205, 4).Value = 44663
$ws.Cells.Item(205, 9).Value = "Primera"
$ws.Cells.Item(205, 10).Value = 200
$ws.Cells.Item(205, 11).Value = 600
$ws.Cells.Item(205, 12).Value = 700
$ws.Cells.Item(205, 13).Value = 650
$ws.Cells.Item(205, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(205, 15).Value = "Región de Ñuble"
$ws.Cells.Item(205, 16).Value = 650
$ws.Cells.Item(205, 17).Value = 1

$ws.Cells.Item(206, 4).Value = 44663
$ws.Cells.Item(206, 9).Value = "Segunda"
$ws.Cells.Item(206, 10).Value = 100
$ws.Cells.Item(206, 11).Value = 500
$ws.Cells.Item(206, 12).Value = 500
$ws.Cells.Item(206, 13).Value = 500
$ws.Cells.Item(206, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(206, 15).Value = "Región de Ñuble"
$ws.Cells.Item(206, 16).Value = 500
$ws.Cells.Item(206, 17).Value = 1

$ws.Cells.Item(207, 4).Value = 44490
$ws.Cells.Item(207, 9).Value = "Primera"
$ws.Cells.Item(207, 10).Value = 200
$ws.Cells.Item(207, 11).Value = 600
$ws.Cells.Item(207, 12).Value = 700
$ws.Cells.Item(207, 13).Value = 650
$ws.Cells.Item(207, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(207, 15).Value = "Región de Ñuble"
$ws.Cells.Item(207, 16).Value = 650
$ws.Cells.Item(207, 17).Value = 1

$ws.Cells.Item(208, 4).Value = 44490
$ws.Cells.Item(208, 9).Value = "Segunda"
$ws.Cells.Item(208, 10).Value = 100
$ws.Cells.Item(208, 11).Value = 500
$ws.Cells.Item(208, 12).Value = 500
$ws.Cells.Item(208, 13).Value = 500
$ws.Cells.Item(208, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(208, 15).Value = "Región de Ñuble"
$ws.Cells.Item(208, 16).Value = 500
$ws.Cells.Item(208, 17).Value = 1

$ws.Cells.Item(209, 4).Value = 44237
$ws.Cells.Item(209, 9).Value = "Primera"
$ws.Cells.Item(209, 10).Value = 200
$ws.Cells.Item(209, 11).Value = 600
$ws.Cells.Item(209, 12).Value = 700
$ws.Cells.Item(209, 13).Value = 650
$ws.Cells.Item(209, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(209, 15).Value = "Región de Ñuble"
$ws.Cells.Item(209, 16).Value = 650
$ws.Cells.Item(209, 17).Value = 1

$ws.Cells.Item(210, 4).Value = 44237
$ws.Cells.Item(210, 9).Value = "Segunda"
$ws.Cells.Item(210, 10).Value = 100
$ws.Cells.Item(210, 11).Value = 500
$ws.Cells.Item(210, 12).Value = 500
$ws.Cells.Item(210, 13).Value = 500
$ws.Cells.Item(210, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(210, 15).Value = "Región de Ñuble"
$ws.Cells.Item(210, 16).Value = 500
$ws.Cells.Item(210, 17).Value = 1

$ws.Cells.Item(211, 4).Value = 44740
$ws.Cells.Item(211, 9).Value = "Primera"
$ws.Cells.Item(211, 10).Value = 100
$ws.Cells.Item(211, 11).Value = 11000
$ws.Cells.Item(211, 12).Value = 12000
$ws.Cells.Item(211, 13).Value = 11500
$ws.Cells.Item(211, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(211, 15).Value = "Región Metropolitana"
$ws.Cells.Item(211, 16).Value = 319
$ws.Cells.Item(211, 17).Value = 36

$ws.Cells.Item(212, 4).Value = 44895
$ws.Cells.Item(212, 9).Value = "Primera"
$ws.Cells.Item(212, 10).Value = 200
$ws.Cells.Item(212, 11).Value = 700
$ws.Cells.Item(212, 12).Value = 800
$ws.Cells.Item(212, 13).Value = 750
$ws.Cells.Item(212, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(212, 15).Value = "Región de Ñuble"
$ws.Cells.Item(212, 16).Value = 750
$ws.Cells.Item(212, 17).Value = 1

$ws.Cells.Item(213, 4).Value = 44895
$ws.Cells.Item(213, 9).Value = "Segunda"
$ws.Cells.Item(213, 10).Value = 100
$ws.Cells.Item(213, 11).Value = 600
$ws.Cells.Item(213, 12).Value = 600
$ws.Cells.Item(213, 13).Value = 600
$ws.Cells.Item(213, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(213, 15).Value = "Región de Ñuble"
$ws.Cells.Item(213, 16).Value = 600
$ws.Cells.Item(213, 17).Value = 1

$ws.Cells.Item(214, 4).Value = 44386
$ws.Cells.Item(214, 9).Value = "Primera"
$ws.Cells.Item(214, 10).Value = 200
$ws.Cells.Item(214, 11).Value = 600
$ws.Cells.Item(214, 12).Value = 700
$ws.Cells.Item(214, 13).Value = 650
$ws.Cells.Item(214, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(214, 15).Value = "Región de Ñuble"
$ws.Cells.Item(214, 16).Value = 650
$ws.Cells.Item(214, 17).Value = 1

$ws.Cells.Item(215, 4).Value = 44386
$ws.Cells.Item(215, 9).Value = "Segunda"
$ws.Cells.Item(215, 10).Value = 100
$ws.Cells.Item(215, 11).Value = 500
$ws.Cells.Item(215, 12).Value = 500
$ws.Cells.Item(215, 13).Value = 500
$ws.Cells.Item(215, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(215, 15).Value = "Región de Ñuble"
$ws.Cells.Item(215, 16).Value = 500
$ws.Cells.Item(215, 17).Value = 1

$ws.Cells.Item(216, 4).Value = 44425
$ws.Cells.Item(216, 9).Value = "Primera"
$ws.Cells.Item(216, 10).Value = 200
$ws.Cells.Item(216, 11).Value = 600
$ws.Cells.Item(216, 12).Value = 700
$ws.Cells.Item(216, 13).Value = 650
$ws.Cells.Item(216, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(216, 15).Value = "Región de Ñuble"
$ws.Cells.Item(216, 16).Value = 650
$ws.Cells.Item(216, 17).Value = 1

$ws.Cells.Item(217, 4).Value = 44425
$ws.Cells.Item(217, 9).Value = "Segunda"
$ws.Cells.Item(217, 10).Value = 100
$ws.Cells.Item(217, 11).Value = 500
$ws.Cells.Item(217, 12).Value = 500
$ws.Cells.Item(217, 13).Value = 500
$ws.Cells.Item(217, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(217, 15).Value = "Región de Ñuble"
$ws.Cells.Item(217, 16).Value = 500
$ws.Cells.Item(217, 17).Value = 1

$ws.Cells.Item(218, 4).Value = 44656
$ws.Cells.Item(218, 9).Value = "Primera"
$ws.Cells.Item(218, 10).Value = 130
$ws.Cells.Item(218, 11).Value = 650
$ws.Cells.Item(218, 12).Value = 6000
$ws.Cells.Item(218, 13).Value = 3942
$ws.Cells.Item(218, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(218, 15).Value = "Región Metropolitana"
$ws.Cells.Item(218, 16).Value = 110
$ws.Cells.Item(218, 17).Value = 36

$ws.Cells.Item(219, 4).Value = 44292
$ws.Cells.Item(219, 9).Value = "Primera"
$ws.Cells.Item(219, 10).Value = 200
$ws.Cells.Item(219, 11).Value = 600
$ws.Cells.Item(219, 12).Value = 700
$ws.Cells.Item(219, 13).Value = 650
$ws.Cells.Item(219, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(219, 15).Value = "Región de Ñuble"
$ws.Cells.Item(219, 16).Value = 650
$ws.Cells.Item(219, 17).Value = 1

$ws.Cells.Item(220, 4).Value = 44292
$ws.Cells.Item(220, 9).Value = "Segunda"
$ws.Cells.Item(220, 10).Value = 100
$ws.Cells.Item(220, 11).Value = 500
$ws.Cells.Item(220, 12).Value = 500
$ws.Cells.Item(220, 13).Value = 500
$ws.Cells.Item(220, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(220, 15).Value = "Región de Ñuble"
$ws.Cells.Item(220, 16).Value = 500
$ws.Cells.Item(220, 17).Value = 1

$ws.Cells.Item(221, 4).Value = 44714
$ws.Cells.Item(221, 9).Value = "Primera"
$ws.Cells.Item(221, 10).Value = 160
$ws.Cells.Item(221, 11).Value = 5000
$ws.Cells.Item(221, 12).Value = 5500
$ws.Cells.Item(221, 13).Value = 5250
$ws.Cells.Item(221, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(221, 15).Value = "Región Metropolitana"
$ws.Cells.Item(221, 16).Value = 146
$ws.Cells.Item(221, 17).Value = 36

$ws.Cells.Item(222, 4).Value = 44264
$ws.Cells.Item(222, 9).Value = "Primera"
$ws.Cells.Item(222, 10).Value = 200
$ws.Cells.Item(222, 11).Value = 600
$ws.Cells.Item(222, 12).Value = 700
$ws.Cells.Item(222, 13).Value = 650
$ws.Cells.Item(222, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(222, 15).Value = "Región de Ñuble"
$ws.Cells.Item(222, 16).Value = 650
$ws.Cells.Item(222, 17).Value = 1

$ws.Cells.Item(223, 4).Value = 44264
$ws.Cells.Item(223, 9).Value = "Segunda"
$ws.Cells.Item(223, 10).Value = 100
$ws.Cells.Item(223, 11).Value = 500
$ws.Cells.Item(223, 12).Value = 500
$ws.Cells.Item(223, 13).Value = 500
$ws.Cells.Item(223, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(223, 15).Value = "Región de Ñuble"
$ws.Cells.Item(223, 16).Value = 500
$ws.Cells.Item(223, 17).Value = 1

$ws.Cells.Item(224, 4).Value = 44376
$ws.Cells.Item(224, 9).Value = "Primera"
$ws.Cells.Item(224, 10).Value = 200
$ws.Cells.Item(224, 11).Value = 600
$ws.Cells.Item(224, 12).Value = 700
$ws.Cells.Item(224, 13).Value = 650
$ws.Cells.Item(224, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(224, 15).Value = "Región de Ñuble"
$ws.Cells.Item(224, 16).Value = 650
$ws.Cells.Item(224, 17).Value = 1

$ws.Cells.Item(225, 4).Value = 44376
$ws.Cells.Item(225, 9).Value = "Segunda"
$ws.Cells.Item(225, 10).Value = 100
$ws.Cells.Item(225, 11).Value = 500
$ws.Cells.Item(225, 12).Value = 500
$ws.Cells.Item(225, 13).Value = 500
$ws.Cells.Item(225, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(225, 15).Value = "Región de Ñuble"
$ws.Cells.Item(225, 16).Value = 500
$ws.Cells.Item(225, 17).Value = 1

$ws.Cells.Item(226, 4).Value = 44847
$ws.Cells.Item(226, 9).Value = "Primera"
$ws.Cells.Item(226, 10).Value = 130
$ws.Cells.Item(226, 11).Value = 5000
$ws.Cells.Item(226, 12).Value = 5500
$ws.Cells.Item(226, 13).Value = 5192
$ws.Cells.Item(226, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(226, 15).Value = "Región Metropolitana"
$ws.Cells.Item(226, 16).Value = 144
$ws.Cells.Item(226, 17).Value = 36

$ws.Cells.Item(227, 4).Value = 44839
$ws.Cells.Item(227, 9).Value = "Primera"
$ws.Cells.Item(227, 10).Value = 200
$ws.Cells.Item(227, 11).Value = 700
$ws.Cells.Item(227, 12).Value = 800
$ws.Cells.Item(227, 13).Value = 750
$ws.Cells.Item(227, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(227, 15).Value = "Región de Ñuble"
$ws.Cells.Item(227, 16).Value = 750
$ws.Cells.Item(227, 17).Value = 1

$ws.Cells.Item(228, 4).Value = 44839
$ws.Cells.Item(228, 9).Value = "Segunda"
$ws.Cells.Item(228, 10).Value = 100
$ws.Cells.Item(228, 11).Value = 600
$ws.Cells.Item(228, 12).Value = 600
$ws.Cells.Item(228, 13).Value = 600
$ws.Cells.Item(228, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(228, 15).Value = "Región de Ñuble"
$ws.Cells.Item(228, 16).Value = 600
$ws.Cells.Item(228, 17).Value = 1

$ws.Cells.Item(229, 4).Value = 44299
$ws.Cells.Item(229, 9).Value = "Primera"
$ws.Cells.Item(229, 10).Value = 200
$ws.Cells.Item(229, 11).Value = 600
$ws.Cells.Item(229, 12).Value = 700
$ws.Cells.Item(229, 13).Value = 650
$ws.Cells.Item(229, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(229, 15).Value = "Región de Ñuble"
$ws.Cells.Item(229, 16).Value = 650
$ws.Cells.Item(229, 17).Value = 1

$ws.Cells.Item(230, 4).Value = 44299
$ws.Cells.Item(230, 9).Value = "Segunda"
$ws.Cells.Item(230, 10).Value = 100
$ws.Cells.Item(230, 11).Value = 500
$ws.Cells.Item(230, 12).Value = 500
$ws.Cells.Item(230, 13).Value = 500
$ws.Cells.Item(230, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(230, 15).Value = "Región de Ñuble"
$ws.Cells.Item(230, 16).Value = 500
$ws.Cells.Item(230, 17).Value = 1

$ws.Cells.Item(231, 4).Value = 44756
$ws.Cells.Item(231, 9).Value = "Primera"
$ws.Cells.Item(231, 10).Value = 200
$ws.Cells.Item(231, 11).Value = 700
$ws.Cells.Item(231, 12).Value = 800
$ws.Cells.Item(231, 13).Value = 750
$ws.Cells.Item(231, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(231, 15).Value = "Región de Ñuble"
$ws.Cells.Item(231, 16).Value = 750
$ws.Cells.Item(231, 17).Value = 1

$ws.Cells.Item(232, 4).Value = 44756
$ws.Cells.Item(232, 9).Value = "Segunda"
$ws.Cells.Item(232, 10).Value = 100
$ws.Cells.Item(232, 11).Value = 600
$ws.Cells.Item(232, 12).Value = 600
$ws.Cells.Item(232, 13).Value = 600
$ws.Cells.Item(232, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(232, 15).Value = "Región de Ñuble"
$ws.Cells.Item(232, 16).Value = 600
$ws.Cells.Item(232, 17).Value = 1

$ws.Cells.Item(233, 4).Value = 44818
$ws.Cells.Item(233, 9).Value = "Primera"
$ws.Cells.Item(233, 10).Value = 90
$ws.Cells.Item(233, 11).Value = 5000
$ws.Cells.Item(233, 12).Value = 6000
$ws.Cells.Item(233, 13).Value = 5556
$ws.Cells.Item(233, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(233, 15).Value = "Región Metropolitana"
$ws.Cells.Item(233, 16).Value = 154
$ws.Cells.Item(233, 17).Value = 36

$ws.Cells.Item(234, 4).Value = 44453
$ws.Cells.Item(234, 9).Value = "Primera"
$ws.Cells.Item(234, 10).Value = 200
$ws.Cells.Item(234, 11).Value = 600
$ws.Cells.Item(234, 12).Value = 700
$ws.Cells.Item(234, 13).Value = 650
$ws.Cells.Item(234, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(234, 15).Value = "Región de Ñuble"
$ws.Cells.Item(234, 16).Value = 650
$ws.Cells.Item(234, 17).Value = 1

$ws.Cells.Item(235, 4).Value = 44453
$ws.Cells.Item(235, 9).Value = "Segunda"
$ws.Cells.Item(235, 10).Value = 100
$ws.Cells.Item(235, 11).Value = 500
$ws.Cells.Item(235, 12).Value = 500
$ws.Cells.Item(235, 13).Value = 500
$ws.Cells.Item(235, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(235, 15).Value = "Región de Ñuble"
$ws.Cells.Item(235, 16).Value = 500
$ws.Cells.Item(235, 17).Value = 1

$ws.Cells.Item(236, 4).Value = 44687
$ws.Cells.Item(236, 9).Value = "Primera"
$ws.Cells.Item(236, 10).Value = 170
$ws.Cells.Item(236, 11).Value = 6500
$ws.Cells.Item(236, 12).Value = 7000
$ws.Cells.Item(236, 13).Value = 6765
$ws.Cells.Item(236, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(236, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(236, 16).Value = 188
$ws.Cells.Item(236, 17).Value = 36

$ws.Cells.Item(237, 4).Value = 44665
$ws.Cells.Item(237, 9).Value = "Primera"
$ws.Cells.Item(237, 10).Value = 200
$ws.Cells.Item(237, 11).Value = 600
$ws.Cells.Item(237, 12).Value = 700
$ws.Cells.Item(237, 13).Value = 650
$ws.Cells.Item(237, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(237, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(237, 16).Value = 650
$ws.Cells.Item(237, 17).Value = 1

$ws.Cells.Item(238, 4).Value = 44665
$ws.Cells.Item(238, 9).Value = "Segunda"
$ws.Cells.Item(238, 10).Value = 100
$ws.Cells.Item(238, 11).Value = 500
$ws.Cells.Item(238, 12).Value = 500
$ws.Cells.Item(238, 13).Value = 500
$ws.Cells.Item(238, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(238, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(238, 16).Value = 500
$ws.Cells.Item(238, 17).Value = 1

$ws.Cells.Item(239, 4).Value = 44813
$ws.Cells.Item(239, 9).Value = "Primera"
$ws.Cells.Item(239, 10).Value = 200
$ws.Cells.Item(239, 11).Value = 5000
$ws.Cells.Item(239, 12).Value = 5500
$ws.Cells.Item(239, 13).Value = 5250
$ws.Cells.Item(239, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(239, 15).Value = "Región Metropolitana"
$ws.Cells.Item(239, 16).Value = 146
$ws.Cells.Item(239, 17).Value = 36

$ws.Cells.Item(240, 4).Value = 44217
$ws.Cells.Item(240, 9).Value = "Primera"
$ws.Cells.Item(240, 10).Value = 200
$ws.Cells.Item(240, 11).Value = 600
$ws.Cells.Item(240, 12).Value = 700
$ws.Cells.Item(240, 13).Value = 650
$ws.Cells.Item(240, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(240, 15).Value = "Región de Ñuble"
$ws.Cells.Item(240, 16).Value = 650
$ws.Cells.Item(240, 17).Value = 1

$ws.Cells.Item(241, 4).Value = 44217
$ws.Cells.Item(241, 9).Value = "Segunda"
$ws.Cells.Item(241, 10).Value = 100
$ws.Cells.Item(241, 11).Value = 500
$ws.Cells.Item(241, 12).Value = 500
$ws.Cells.Item(241, 13).Value = 500
$ws.Cells.Item(241, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(241, 15).Value = "Región de Ñuble"
$ws.Cells.Item(241, 16).Value = 500
$ws.Cells.Item(241, 17).Value = 1

$ws.Cells.Item(242, 4).Value = 44679
$ws.Cells.Item(242, 9).Value = "Primera"
$ws.Cells.Item(242, 10).Value = 200
$ws.Cells.Item(242, 11).Value = 600
$ws.Cells.Item(242, 12).Value = 700
$ws.Cells.Item(242, 13).Value = 650
$ws.Cells.Item(242, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(242, 15).Value = "Región de Ñuble"
$ws.Cells.Item(242, 16).Value = 650
$ws.Cells.Item(242, 17).Value = 1

$ws.Cells.Item(243, 4).Value = 44679
$ws.Cells.Item(243, 9).Value = "Segunda"
$ws.Cells.Item(243, 10).Value = 100
$ws.Cells.Item(243, 11).Value = 500
$ws.Cells.Item(243, 12).Value = 500
$ws.Cells.Item(243, 13).Value = 500
$ws.Cells.Item(243, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(243, 15).Value = "Región de Ñuble"
$ws.Cells.Item(243, 16).Value = 500
$ws.Cells.Item(243, 17).Value = 1

$ws.Cells.Item(244, 4).Value = 44350
$ws.Cells.Item(244, 9).Value = "Primera"
$ws.Cells.Item(244, 10).Value = 200
$ws.Cells.Item(244, 11).Value = 600
$ws.Cells.Item(244, 12).Value = 700
$ws.Cells.Item(244, 13).Value = 650
$ws.Cells.Item(244, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(244, 15).Value = "Región de Ñuble"
$ws.Cells.Item(244, 16).Value = 650
$ws.Cells.Item(244, 17).Value = 1

$ws.Cells.Item(245, 4).Value = 44350
$ws.Cells.Item(245, 9).Value = "Segunda"
$ws.Cells.Item(245, 10).Value = 100
$ws.Cells.Item(245, 11).Value = 500
$ws.Cells.Item(245, 12).Value = 500
$ws.Cells.Item(245, 13).Value = 500
$ws.Cells.Item(245, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(245, 15).Value = "Región de Ñuble"
$ws.Cells.Item(245, 16).Value = 500
$ws.Cells.Item(245, 17).Value = 1

$ws.Cells.Item(246, 4).Value = 44890
$ws.Cells.Item(246, 9).Value = "Primera"
$ws.Cells.Item(246, 10).Value = 270
$ws.Cells.Item(246, 11).Value = 15000
$ws.Cells.Item(246, 12).Value = 16000
$ws.Cells.Item(246, 13).Value = 15444
$ws.Cells.Item(246, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(246, 15).Value = "Región Metropolitana"
$ws.Cells.Item(246, 16).Value = 429
$ws.Cells.Item(246, 17).Value = 36

$ws.Cells.Item(247, 4).Value = 44447
$ws.Cells.Item(247, 9).Value = "Primera"
$ws.Cells.Item(247, 10).Value = 200
$ws.Cells.Item(247, 11).Value = 600
$ws.Cells.Item(247, 12).Value = 700
$ws.Cells.Item(247, 13).Value = 650
$ws.Cells.Item(247, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(247, 15).Value = "Región de Ñuble"
$ws.Cells.Item(247, 16).Value = 650
$ws.Cells.Item(247, 17).Value = 1

$ws.Cells.Item(248, 4).Value = 44447
$ws.Cells.Item(248, 9).Value = "Segunda"
$ws.Cells.Item(248, 10).Value = 100
$ws.Cells.Item(248, 11).Value = 500
$ws.Cells.Item(248, 12).Value = 500
$ws.Cells.Item(248, 13).Value = 500
$ws.Cells.Item(248, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(248, 15).Value = "Región de Ñuble"
$ws.Cells.Item(248, 16).Value = 500
$ws.Cells.Item(248, 17).Value = 1

$ws.Cells.Item(249, 4).Value = 44763
$ws.Cells.Item(249, 9).Value = "Primera"
$ws.Cells.Item(249, 10).Value = 35
$ws.Cells.Item(249, 11).Value = 17000
$ws.Cells.Item(249, 12).Value = 18000
$ws.Cells.Item(249, 13).Value = 17429
$ws.Cells.Item(249, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(249, 15).Value = "Región Metropolitana"
$ws.Cells.Item(249, 16).Value = 484
$ws.Cells.Item(249, 17).Value = 36

$ws.Cells.Item(250, 4).Value = 44245
$ws.Cells.Item(250, 9).Value = "Primera"
$ws.Cells.Item(250, 10).Value = 200
$ws.Cells.Item(250, 11).Value = 600
$ws.Cells.Item(250, 12).Value = 700
$ws.Cells.Item(250, 13).Value = 650
$ws.Cells.Item(250, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(250, 15).Value = "Región de Ñuble"
$ws.Cells.Item(250, 16).Value = 650
$ws.Cells.Item(250, 17).Value = 1

$ws.Cells.Item(251, 4).Value = 44245
$ws.Cells.Item(251, 9).Value = "Segunda"
$ws.Cells.Item(251, 10).Value = 100
$ws.Cells.Item(251, 11).Value = 500
$ws.Cells.Item(251, 12).Value = 500
$ws.Cells.Item(251, 13).Value = 500
$ws.Cells.Item(251, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(251, 15).Value = "Región de Ñuble"
$ws.Cells.Item(251, 16).Value = 500
$ws.Cells.Item(251, 17).Value = 1

$ws.Cells.Item(252, 4).Value = 44565
$ws.Cells.Item(252, 9).Value = "Primera"
$ws.Cells.Item(252, 10).Value = 200
$ws.Cells.Item(252, 11).Value = 600
$ws.Cells.Item(252, 12).Value = 700
$ws.Cells.Item(252, 13).Value = 650
$ws.Cells.Item(252, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(252, 15).Value = "Región de Ñuble"
$ws.Cells.Item(252, 16).Value = 650
$ws.Cells.Item(252, 17).Value = 1

$ws.Cells.Item(253, 4).Value = 44565
$ws.Cells.Item(253, 9).Value = "Segunda"
$ws.Cells.Item(253, 10).Value = 100
$ws.Cells.Item(253, 11).Value = 500
$ws.Cells.Item(253, 12).Value = 500
$ws.Cells.Item(253, 13).Value = 500
$ws.Cells.Item(253, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(253, 15).Value = "Región de Ñuble"
$ws.Cells.Item(253, 16).Value = 500
$ws.Cells.Item(253, 17).Value = 1

$ws.Cells.Item(254, 4).Value = 44806
$ws.Cells.Item(254, 9).Value = "Primera"
$ws.Cells.Item(254, 10).Value = 220
$ws.Cells.Item(254, 11).Value = 5000
$ws.Cells.Item(254, 12).Value = 5500
$ws.Cells.Item(254, 13).Value = 5273
$ws.Cells.Item(254, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(254, 15).Value = "Región Metropolitana"
$ws.Cells.Item(254, 16).Value = 146
$ws.Cells.Item(254, 17).Value = 36

$ws.Cells.Item(255, 4).Value = 44911
$ws.Cells.Item(255, 9).Value = "Primera"
$ws.Cells.Item(255, 10).Value = 200
$ws.Cells.Item(255, 11).Value = 700
$ws.Cells.Item(255, 12).Value = 800
$ws.Cells.Item(255, 13).Value = 750
$ws.Cells.Item(255, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(255, 15).Value = "Región de Ñuble"
$ws.Cells.Item(255, 16).Value = 750
$ws.Cells.Item(255, 17).Value = 1

$ws.Cells.Item(256, 4).Value = 44911
$ws.Cells.Item(256, 9).Value = "Segunda"
$ws.Cells.Item(256, 10).Value = 100
$ws.Cells.Item(256, 11).Value = 600
$ws.Cells.Item(256, 12).Value = 600
$ws.Cells.Item(256, 13).Value = 600
$ws.Cells.Item(256, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(256, 15).Value = "Región de Ñuble"
$ws.Cells.Item(256, 16).Value = 600
$ws.Cells.Item(256, 17).Value = 1

$ws.Cells.Item(257, 1).Value = 11
$ws.Cells.Item(257, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(257, 3).Value = "Bíobío"
$ws.Cells.Item(257, 4).Value = 44736
$ws.Cells.Item(257, 4).NumberFormat = $dateFormat
$ws.Cells.Item(257, 5).Value = 8
$ws.Cells.Item(257, 6).Value = 100112040
$ws.Cells.Item(257, 7).Value = "Cilantro"
$ws.Cells.Item(257, 8).Value = "Sin especificar"
$ws.Cells.Item(257, 9).Value = "Primera"
$ws.Cells.Item(257, 10).Value = 200
$ws.Cells.Item(257, 11).Value = 600
$ws.Cells.Item(257, 12).Value = 700
$ws.Cells.Item(257, 13).Value = 650
$ws.Cells.Item(257, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(257, 15).Value = "Región de Ñuble"
$ws.Cells.Item(257, 16).Value = 650
$ws.Cells.Item(257, 17).Value = 1
$ws.Cells.Item(257, 18).Value = "Hortaliza"

$ws.Cells.Item(258, 1).Value = 11
$ws.Cells.Item(258, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(258, 3).Value = "Bíobío"
$ws.Cells.Item(258, 4).Value = 44736
$ws.Cells.Item(258, 4).NumberFormat = $dateFormat
$ws.Cells.Item(258, 5).Value = 8
$ws.Cells.Item(258, 6).Value = 100112040
$ws.Cells.Item(258, 7).Value = "Cilantro"
$ws.Cells.Item(258, 8).Value = "Sin especificar"
$ws.Cells.Item(258, 9).Value = "Segunda"
$ws.Cells.Item(258, 10).Value = 100
$ws.Cells.Item(258, 11).Value = 500
$ws.Cells.Item(258, 12).Value = 500
$ws.Cells.Item(258, 13).Value = 500
$ws.Cells.Item(258, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(258, 15).Value = "Región de Ñuble"
$ws.Cells.Item(258, 16).Value = 500
$ws.Cells.Item(258, 17).Value = 1
$ws.Cells.Item(258, 18).Value = "Hortaliza"
